# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff/Handback
# Datetime" timestamps to reflect a newer report generation run.

$wb = $excel.ActiveWorkbook

# "Latest HO Xliff Generate Date" (Overview) and "Correspond Handoff
# Datetime" (de-de) for 8cc863a9-...md share the same underlying value,
# so both move to the new handoff-generation timestamp.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-20 23:08:06"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-20 23:08:06"
$dede.Range("K2").Value = "2016-08-20 23:08:33"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-20 23:07:58"
$zhcn.Range("K2").Value = "2016-08-20 23:08:25"
